$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.726598034398364
$ws.Range("D2").Value = 0.0202251898184258
$ws.Range("E2").Value = 1.323723653710346
$ws.Range("F2").Value = 0.2722498823717743
$ws.Range("G2").Value = 0.1412063102096681
$ws.Range("H2").Value = 0.2988846374982899
$ws.Range("L2").Value = 0.8459651622079605
$ws.Range("O2").Value = 0.7769765906976716
$ws.Range("B3").Value = 1.592131924085265
$ws.Range("D3").Value = 0.01764334772118303
$ws.Range("E3").Value = 1.217921524984206
$ws.Range("F3").Value = 0.270144820187241
$ws.Range("G3").Value = 0.1396760083086903
$ws.Range("H3").Value = 0.3029791293159576
$ws.Range("L3").Value = 0.7414064603125325
$ws.Range("O3").Value = 0.7819990035903857
$ws.Range("B4").Value = 1.509456002926925
$ws.Range("D4").Value = 0.01605168090171105
$ws.Range("E4").Value = 1.153363124299545
$ws.Range("F4").Value = 0.2692320846979968
$ws.Range("G4").Value = 0.139029638160892
$ws.Range("H4").Value = 0.3057860215178323
$ws.Range("L4").Value = 0.6769404639647689
$ws.Range("O4").Value = 0.786307230825912
$ws.Range("B5").Value = 1.47573751472811
$ws.Range("D5").Value = 0.01540151803177991
$ws.Range("E5").Value = 1.127160117646127
$ws.Range("F5").Value = 0.2689550014387336
$ws.Range("G5").Value = 0.1388392021241813
$ws.Range("H5").Value = 0.3070032221731722
$ws.Range("L5").Value = 0.6506042941818464
$ws.Range("O5").Value = 0.7883685505049982
$ws.Range("B6").Value = 1.470136965125334
$ws.Range("D6").Value = 0.01529346755139471
$ws.Range("E6").Value = 1.12281557416361
$ws.Range("F6").Value = 0.2689147028416272
$ws.Range("G6").Value = 0.1388119645618389
$ws.Range("H6").Value = 0.3072097607365833
$ws.Range("L6").Value = 0.6462272605396606
$ws.Range("O6").Value = 0.7887292284226248
$ws.Range("B7").Value = 1.509001372989871
$ws.Range("D7").Value = 0.01604291874945574
$ws.Range("E7").Value = 1.153009311487295
$ws.Range("F7").Value = 0.2692279646014129
$ws.Range("G7").Value = 0.1390267754679897
$ws.Range("H7").Value = 0.3058021404057527
$ws.Range("L7").Value = 0.6765855496630024
$ws.Range("O7").Value = 0.7863337957557235
$ws.Range("B8").Value = 1.680258198473325
$ws.Range("D8").Value = 0.01933633800719292
$ws.Range("E8").Value = 1.287161058469138
$ws.Range("F8").Value = 0.2714448508545217
$ws.Range("G8").Value = 0.1406173353320668
$ws.Range("H8").Value = 0.3002354725605656
$ws.Range("L8").Value = 0.8099693508221435
$ws.Range("O8").Value = 0.7784528548934304
$ws.Range("B9").Value = 2.01516667354332
$ws.Range("D9").Value = 0.02574138222941258
$ws.Range("E9").Value = 1.553313630393262
$ws.Range("F9").Value = 0.2788346860656574
$ws.Range("G9").Value = 0.1460981166286714
$ws.Range("H9").Value = 0.2916551389581628
$ws.Range("L9").Value = 1.069379513823549
$ws.Range("O9").Value = 0.7728137080831772
$ws.Range("B10").Value = 2.260652077066879
$ws.Range("D10").Value = 0.030411759826805
$ws.Range("E10").Value = 1.750588907990704
$ws.Range("F10").Value = 0.2861598270207537
$ws.Range("G10").Value = 0.1516131263038289
$ws.Range("H10").Value = 0.2867919412197324
$ws.Range("L10").Value = 1.258623019224615
$ws.Range("O10").Value = 0.7747926956490687
$ws.Range("B11").Value = 2.372206965000771
$ws.Range("D11").Value = 0.0325281777344486
$ws.Range("E11").Value = 1.840683230460826
$ws.Range("F11").Value = 0.2899126815779027
$ws.Range("G11").Value = 0.1544555883509418
$ws.Range("H11").Value = 0.2848959747797153
$ws.Range("L11").Value = 1.34441850880421
$ws.Range("O11").Value = 0.7770519127064972
$ws.Range("B12").Value = 2.414432624694882
$ws.Range("D12").Value = 0.03332838044360642
$ws.Range("E12").Value = 1.874847721823613
$ws.Range("F12").Value = 0.2913949588415861
$ws.Range("G12").Value = 0.1555807584664421
$ws.Range("H12").Value = 0.2842238000705777
$ws.Range("L12").Value = 1.376864304998662
$ws.Range("O12").Value = 0.7781052005724973
$ws.Range("B13").Value = 2.405339375601386
$ws.Range("D13").Value = 0.03315609860072755
$ws.Range("E13").Value = 1.867487711478162
$ws.Range("F13").Value = 0.2910729935108023
$ws.Range("G13").Value = 0.155336250044499
$ws.Range("H13").Value = 0.2843665237962796
$ws.Range("L13").Value = 1.369878455006756
$ws.Range("O13").Value = 0.7778695228732886
$ws.Range("B14").Value = 2.37568125582078
$ws.Range("D14").Value = 0.03259403609168032
$ws.Range("E14").Value = 1.843493024078839
$ws.Range("F14").Value = 0.2900334000089728
$ws.Range("G14").Value = 0.1545471742813049
$ws.Range("H14").Value = 0.2848397552200908
$ws.Range("L14").Value = 1.347088713469645
$ws.Range("O14").Value = 0.7771345902530982
$ws.Range("B15").Value = 2.357512465206696
$ws.Range("D15").Value = 0.03224959332129629
$ws.Range("E15").Value = 1.828801713117457
$ws.Range("F15").Value = 0.2894046026574415
$ws.Range("G15").Value = 0.1540702207183386
$ws.Range("H15").Value = 0.2851355947621528
$ws.Range("L15").Value = 1.333123696404698
$ws.Range("O15").Value = 0.7767102481628285
$ws.Range("B16").Value = 2.253359262129834
$ws.Range("D16").Value = 0.03027327712379702
$ws.Range("E16").Value = 1.744707872237456
$ws.Range("F16").Value = 0.2859230902228234
$ws.Range("G16").Value = 0.1514341482392894
$ws.Range("H16").Value = 0.2869222379959808
$ws.Range("L16").Value = 1.253010107943624
$ws.Range("O16").Value = 0.7746725978036295
$ws.Range("B17").Value = 2.189433893614364
$ws.Range("D17").Value = 0.02905873428153427
$ws.Range("E17").Value = 1.693207216298845
$ws.Range("F17").Value = 0.2838955011577724
$ws.Range("G17").Value = 0.149903039468299
$ws.Range("H17").Value = 0.2880995179237118
$ws.Range("L17").Value = 1.203787370835698
$ws.Range("O17").Value = 0.7737723070120808
$ws.Range("B18").Value = 2.152654683756225
$ws.Range("D18").Value = 0.0283593971258469
$ws.Range("E18").Value = 1.663618782542244
$ws.Range("F18").Value = 0.2827688357016882
$ws.Range("G18").Value = 0.1490537432091514
$ws.Range("H18").Value = 0.2888064145461371
$ws.Range("L18").Value = 1.175448297937635
$ws.Range("O18").Value = 0.7733822557772498
$ws.Range("B19").Value = 2.140200007750707
$ws.Range("D19").Value = 0.02812248433044573
$ws.Range("E19").Value = 1.653606473825079
$ws.Range("F19").Value = 0.2823941384583577
$ws.Range("G19").Value = 0.148771546885925
$ws.Range("H19").Value = 0.2890508591402252
$ws.Range("L19").Value = 1.165848492375744
$ws.Range("O19").Value = 0.7732720642529216
$ws.Range("B20").Value = 2.196240005925574
$ws.Range("D20").Value = 0.02918810401342853
$ws.Range("E20").Value = 1.698686116289196
$ws.Range("F20").Value = 0.2841072433226302
$ws.Range("G20").Value = 0.150062777385088
$ws.Range("H20").Value = 0.2879711126216478
$ws.Range("L20").Value = 1.209030064744184
$ws.Range("O20").Value = 0.7738549034730795
$ws.Range("B21").Value = 2.384393052712539
$ws.Range("D21").Value = 0.03275916161288706
$ws.Range("E21").Value = 1.850539573583688
$ws.Range("F21").Value = 0.2903370888384416
$ws.Range("G21").Value = 0.1547776143205368
$ws.Range("H21").Value = 0.2846995106648365
$ws.Range("L21").Value = 1.353783790643661
$ws.Range("O21").Value = 0.7773450716101706
$ws.Range("B22").Value = 2.507258115939976
$ws.Range("D22").Value = 0.0350858037243853
$ws.Range("E22").Value = 1.950061449549366
$ws.Range("F22").Value = 0.2947654056958413
$ws.Range("G22").Value = 0.1581438010361182
$ws.Range("H22").Value = 0.2828283471916109
$ws.Range("L22").Value = 1.448137090570924
$ws.Range("O22").Value = 0.7807799466666268
$ws.Range("B23").Value = 2.441692477908077
$ws.Range("D23").Value = 0.03384471581816229
$ws.Range("E23").Value = 1.896920385940149
$ws.Range("F23").Value = 0.2923690665662591
$ws.Range("G23").Value = 0.1563208801003952
$ws.Range("H23").Value = 0.2838024916842983
$ws.Range("L23").Value = 1.39780231510349
$ws.Range("O23").Value = 0.7788403368479635
$ws.Range("B24").Value = 2.193163048640258
$ws.Range("D24").Value = 0.02912961931808411
$ws.Range("E24").Value = 1.696209043136577
$ws.Range("F24").Value = 0.2840113931879245
$ws.Range("G24").Value = 0.1499904635016662
$ws.Range("H24").Value = 0.2880290710348987
$ws.Range("L24").Value = 1.206659968201393
$ws.Range("O24").Value = 0.773817164515151
$ws.Range("B25").Value = 1.924665212247646
$ws.Range("D25").Value = 0.02401466106498162
$ws.Range("E25").Value = 1.480997856433788
$ws.Range("F25").Value = 0.276505561267669
$ws.Range("G25").Value = 0.1443575545189262
$ws.Range("H25").Value = 0.2937245368297283
$ws.Range("L25").Value = 0.9994363999876441
$ws.Range("O25").Value = 0.7732741032007624
